# Insert a new weekly record at row 741 (shifts existing rows 741-830 down to 742-831).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(741).Insert()

$ws.Cells.Item(741, 1).Value = 6
$ws.Cells.Item(741, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(741, 3).Value = "Metropolitana"
$ws.Cells.Item(741, 4).Value = 45212
$ws.Cells.Item(741, 5).Value = 13
$ws.Cells.Item(741, 6).Value = 100112044
$ws.Cells.Item(741, 7).Value = "Perejil"
$ws.Cells.Item(741, 8).Value = "Sin especificar"
$ws.Cells.Item(741, 9).Value = "Primera"
$ws.Cells.Item(741, 10).Value = 280
$ws.Cells.Item(741, 11).Value = 15000
$ws.Cells.Item(741, 12).Value = 16000
$ws.Cells.Item(741, 13).Value = 15536
$ws.Cells.Item(741, 14).Value = "`$/docena de atados"
$ws.Cells.Item(741, 15).Value = "Región Metropolitana"
$ws.Cells.Item(741, 16).Value = 5179
$ws.Cells.Item(741, 17).Value = 3
$ws.Cells.Item(741, 18).Value = "Hortaliza"
